$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("question_answers")
$ws2 = $wb.Worksheets.Item("outputs")

# question_answers sheet: column B holds text-typed answer values (1-5)
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "1"
$ws1.Range("B3").NumberFormat = "@"
$ws1.Range("B3").Value = "3"
$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = "2"
$ws1.Range("B5").NumberFormat = "@"
$ws1.Range("B5").Value = "5"
$ws1.Range("B6").NumberFormat = "@"
$ws1.Range("B6").Value = "3"
$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "2"
$ws1.Range("B8").NumberFormat = "@"
$ws1.Range("B8").Value = "5"
$ws1.Range("B9").NumberFormat = "@"
$ws1.Range("B9").Value = "4"
$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = "4"
$ws1.Range("B11").NumberFormat = "@"
$ws1.Range("B11").Value = "5"
$ws1.Range("B13").NumberFormat = "@"
$ws1.Range("B13").Value = "4"
$ws1.Range("B15").NumberFormat = "@"
$ws1.Range("B15").Value = "1"
$ws1.Range("B16").NumberFormat = "@"
$ws1.Range("B16").Value = "4"
$ws1.Range("B17").NumberFormat = "@"
$ws1.Range("B17").Value = "5"
$ws1.Range("B18").NumberFormat = "@"
$ws1.Range("B18").Value = "4"
$ws1.Range("B19").NumberFormat = "@"
$ws1.Range("B19").Value = "4"
$ws1.Range("B20").NumberFormat = "@"
$ws1.Range("B20").Value = "2"
$ws1.Range("B21").NumberFormat = "@"
$ws1.Range("B21").Value = "4"
$ws1.Range("B22").NumberFormat = "@"
$ws1.Range("B22").Value = "4"
$ws1.Range("B23").NumberFormat = "@"
$ws1.Range("B23").Value = "3"
$ws1.Range("B24").NumberFormat = "@"
$ws1.Range("B24").Value = "2"
$ws1.Range("B25").NumberFormat = "@"
$ws1.Range("B25").Value = "1"
$ws1.Range("B26").NumberFormat = "@"
$ws1.Range("B26").Value = "4"
$ws1.Range("B27").NumberFormat = "@"
$ws1.Range("B27").Value = "5"
$ws1.Range("B28").NumberFormat = "@"
$ws1.Range("B28").Value = "5"
$ws1.Range("B29").NumberFormat = "@"
$ws1.Range("B29").Value = "3"
$ws1.Range("B30").NumberFormat = "@"
$ws1.Range("B30").Value = "2"
$ws1.Range("B31").NumberFormat = "@"
$ws1.Range("B31").Value = "5"
$ws1.Range("B32").NumberFormat = "@"
$ws1.Range("B32").Value = "4"
$ws1.Range("B33").NumberFormat = "@"
$ws1.Range("B33").Value = "3"
$ws1.Range("B35").NumberFormat = "@"
$ws1.Range("B35").Value = "3"
$ws1.Range("B36").NumberFormat = "@"
$ws1.Range("B36").Value = "3"
$ws1.Range("B37").NumberFormat = "@"
$ws1.Range("B37").Value = "3"
$ws1.Range("B38").NumberFormat = "@"
$ws1.Range("B38").Value = "2"
$ws1.Range("B39").NumberFormat = "@"
$ws1.Range("B39").Value = "3"
$ws1.Range("B40").NumberFormat = "@"
$ws1.Range("B40").Value = "4"
$ws1.Range("B41").NumberFormat = "@"
$ws1.Range("B41").Value = "3"
$ws1.Range("B42").NumberFormat = "@"
$ws1.Range("B42").Value = "3"
$ws1.Range("B43").NumberFormat = "@"
$ws1.Range("B43").Value = "3"
$ws1.Range("B45").NumberFormat = "@"
$ws1.Range("B45").Value = "2"
$ws1.Range("B46").NumberFormat = "@"
$ws1.Range("B46").Value = "3"
$ws1.Range("B47").NumberFormat = "@"
$ws1.Range("B47").Value = "3"
$ws1.Range("B48").NumberFormat = "@"
$ws1.Range("B48").Value = "5"
$ws1.Range("B49").NumberFormat = "@"
$ws1.Range("B49").Value = "1"
$ws1.Range("B50").NumberFormat = "@"
$ws1.Range("B50").Value = "2"
$ws1.Range("B51").NumberFormat = "@"
$ws1.Range("B51").Value = "2"
$ws1.Range("B52").NumberFormat = "@"
$ws1.Range("B52").Value = "4"
$ws1.Range("B53").NumberFormat = "@"
$ws1.Range("B53").Value = "1"

# outputs sheet: column B holds numeric score totals
$ws2.Range("B2").Value = 7
$ws2.Range("B3").Value = 23
$ws2.Range("B4").Value = 17
$ws2.Range("B5").Value = 14
$ws2.Range("B6").Value = 21
$ws2.Range("B7").Value = 27
$ws2.Range("B8").Value = 35
$ws2.Range("B9").Value = 144
